# Adicionando suporte a NotaReferenciada
#
# Adds four new styles to the document's style sheet, mirroring the
# diff applied to reference.docx:
#   - NotaReferenciada   (character, based on DefaultCharacter)
#   - Textodenotadefim   (paragraph, "endnote text", based on Normal)
#   - TextodenotadefimChar (character, "Texto de nota de fim Char", based on Fontepargpadro)
#   - Refdenotadefim     (character, "endnote reference", based on Fontepargpadro, superscript)

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1, wdStyleTypeCharacter = 2

# --- NotaReferenciada ---------------------------------------------------
$notaReferenciada = $d.Styles.Add("NotaReferenciada", 2)
$notaReferenciada.NameLocal = "Nota Referenciada"
$notaReferenciada.BaseStyle = "DefaultCharacter"
$notaReferenciada.Priority = 1
$notaReferenciada.QuickStyle = $true

# --- Textodenotadefim / TextodenotadefimChar (linked pair) --------------
$textoNotaFim = $d.Styles.Add("Textodenotadefim", 1)
$textoNotaFim.NameLocal = "endnote text"
$textoNotaFim.BaseStyle = "Normal"
$textoNotaFim.Priority = 99

$textoNotaFimChar = $d.Styles.Add("TextodenotadefimChar", 2)
$textoNotaFimChar.NameLocal = "Texto de nota de fim Char"
$textoNotaFimChar.BaseStyle = "Fontepargpadro"
$textoNotaFimChar.Priority = 99

$textoNotaFim.LinkStyle = $textoNotaFimChar
$textoNotaFimChar.LinkStyle = $textoNotaFim

# --- Refdenotadefim -------------------------------------------------------
$refNotaFim = $d.Styles.Add("Refdenotadefim", 2)
$refNotaFim.NameLocal = "endnote reference"
$refNotaFim.BaseStyle = "Fontepargpadro"
$refNotaFim.Priority = 99
$refNotaFim.Font.Superscript = $true
